$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host $ws.Name
Write-Host $ws.Range("A2").Value
